# Phase 11 - Onglet 3 Amortissements: update mappings_obligatoires.xlsx
# - Add 4 new "Immobilisation" mapping rows (Facade/Toiture, IGT, agencements, structure/GO)
# - Strip the stray font/cell formatting that had been applied to column A/B/C
# - Restore the selection to A53 and size the (now unused) helper column G
# - Drop calcOnSave from the workbook calculation properties

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the extra cell-level formatting (custom Helvetica / black-colored
#     Aptos Narrow fonts) that used to be applied across A1:C53 so the cells
#     fall back to the workbook's default style. ---
$ws.Range("A1:C53").ClearFormats()

# --- Append the four new "Immobilisations" / "Actif" mapping rows. ---
$newRows = @(
    @("Immobilisation Facade/Toiture", "Immobilisations ", "Actif"),
    @("Immobilisation IGT", "Immobilisations ", "Actif"),
    @("Immobilisation agencements", "Immobilisations ", "Actif"),
    @("Immobilisation structure/GO", "Immobilisations ", "Actif")
)

$startRow = 54
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $rowIndex = $startRow + $i
    $rowData = $newRows[$i]
    $ws.Cells.Item($rowIndex, 1).Value = $rowData[0]
    $ws.Cells.Item($rowIndex, 2).Value = $rowData[1]
    $ws.Cells.Item($rowIndex, 3).Value = $rowData[2]
}

# --- Size the helper column G (leftover custom width from drafting the new
#     mapping labels before they were copied into column A). ---
$ws.Columns("G").ColumnWidth = 27.83

# --- Restore the view: active cell back on the last data row. ---
$ws.Activate()
$ws.Range("A53").Select()

# --- Calculation options: stop forcing a recalculation on every save. ---
$excel.CalculateBeforeSave = $false
